$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12: copy the formatting from row 11's relevant cells (A/C/D/E),
# matching the pattern already used for earlier "Spielen"/"Anna Franziska" rows,
# then fill in this entry's values.

$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)  # xlPasteFormats

# Date: 2017-09-10 (serial 42988)
$ws.Range("A12").Value2 = 42988
$ws.Range("C12").Value = "Spielen"
$ws.Range("D12").Value = "Anna Franziska"
$ws.Range("E12").Value = "Object wird jetzt nur bis zur nächsten Wand geworfen; Poweranzeige ist implementiert; Sichtradius der Wachen korrekt kleiner, wenn Spieler verkleidet"

# Match the row height used for this new, long note (wrapped text).
$ws.Rows.Item(12).RowHeight = 75

# Move the active selection below the newly added row, as in the edited file.
$ws.Range("A13").Select()
